# Apply metadata updates described in the commit "Atualizacoes 16 de janeiro de 2024."
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.1 -> 0.0.0
$ws.Range("B3").Value = "0.0.0"

# Title: Extension of Patient Gender Identity -> Gender Identity
$ws.Range("B5").Value = "Gender Identity"

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-16T20:08:55-03:00
$ws.Range("B8").Value = "2024-01-16T20:08:55-03:00"

# Description: The patient's gender. -> Extension of the patient's gender.
$ws.Range("B12").Value = "Extension of the patient's gender."
